$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update player_id values in column C for rows 2-23 (LeBron James rows)
# from 3462 to 3463, matching the data cleanup commit.
$ws.Range("C2:C23").Value = 3463
